# Auto-generated Excel COM-interop script applying numeric value updates
# to the "Leviathan Profits" workbook (per-sheet leve-profit recalculations).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2958.6667
$ws.Range("I2").Value = 3213.8462
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 3213.8462
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = -3100.8462
$ws.Range("N2").Value = -1526
$ws.Range("H98").Value = 1762.1904
$ws.Range("I98").Value = 1128.8948
$ws.Range("K98").Value = 1128.8948
$ws.Range("M98").Value = 369.1052
$ws.Range("H116").Value = 4332.8
$ws.Range("H122").Value = 1762.1904
$ws.Range("I122").Value = 1128.8948
$ws.Range("K122").Value = 3386.6844
$ws.Range("M122").Value = -936.6844000000001
$ws.Range("H138").Value = 4518
$ws.Range("I138").Value = 4289.4
$ws.Range("J138").Value = 4589.4375
$ws.Range("K138").Value = 12868.2
$ws.Range("L138").Value = 13768.3125
$ws.Range("M138").Value = -7728.199999999999
$ws.Range("N138").Value = -24048.3125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27800
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 27800
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 27800
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = -28374
$ws.Range("H41").Value = 9099.875
$ws.Range("I41").Value = 3999
$ws.Range("J41").Value = 9828.571
$ws.Range("K41").Value = 3999
$ws.Range("L41").Value = 9828.571
$ws.Range("M41").Value = -3585
$ws.Range("N41").Value = -10656.571
$ws.Range("H63").Value = 6913.385
$ws.Range("I63").Value = 4267.7144
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 4267.7144
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -3581.7144
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 6913.385
$ws.Range("I66").Value = 4267.7144
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 21338.572
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -17906.572
$ws.Range("N66").Value = -56864
$ws.Range("H121").Value = 55999
$ws.Range("J121").Value = 55999
$ws.Range("L121").Value = 55999
$ws.Range("N121").Value = -59493
$ws.Range("H122").Value = 2677.625
$ws.Range("J122").Value = 3579.8
$ws.Range("L122").Value = 10739.4
$ws.Range("N122").Value = -15639.4
$ws.Range("H132").Value = 1709.1875
$ws.Range("I132").Value = 1556.711
$ws.Range("K132").Value = 4670.133
$ws.Range("M132").Value = -2140.133

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 679.7
$ws.Range("I64").Value = 707
$ws.Range("K64").Value = 707
$ws.Range("M64").Value = -482
$ws.Range("H67").Value = 679.7
$ws.Range("I67").Value = 707
$ws.Range("K67").Value = 707
$ws.Range("M67").Value = 73

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2033.3334
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -2200
$ws.Range("H31").Value = 3778.25
$ws.Range("I31").Value = 2306.4
$ws.Range("K31").Value = 2306.4
$ws.Range("M31").Value = -2011.4
$ws.Range("H34").Value = 3778.25
$ws.Range("I34").Value = 2306.4
$ws.Range("K34").Value = 2306.4
$ws.Range("M34").Value = -2104.4
$ws.Range("H35").Value = 4269.5713
$ws.Range("I35").Value = 1554.8
$ws.Range("J35").Value = 5777.778
$ws.Range("K35").Value = 1554.8
$ws.Range("L35").Value = 5777.778
$ws.Range("M35").Value = -1260.8
$ws.Range("N35").Value = -6365.778
$ws.Range("H62").Value = 3141.6667
$ws.Range("I62").Value = 3141.6667
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3141.6667
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2517.6667
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 3141.6667
$ws.Range("I65").Value = 3141.6667
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15708.3335
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -12588.3335
$ws.Range("N65").Value = $null
$ws.Range("H122").Value = 68529.8
$ws.Range("I122").Value = 101149.4
$ws.Range("J122").Value = 3290.6
$ws.Range("K122").Value = 303448.2
$ws.Range("L122").Value = 9871.799999999999
$ws.Range("M122").Value = -300998.2
$ws.Range("N122").Value = -14771.8
$ws.Range("H132").Value = 3709.6
$ws.Range("I132").Value = 3709.6
$ws.Range("K132").Value = 11128.8
$ws.Range("M132").Value = -8598.799999999999
$ws.Range("H141").Value = 265817.7
$ws.Range("J141").Value = 265817.7
$ws.Range("L141").Value = 265817.7
$ws.Range("N141").Value = -276177.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1210.5294
$ws.Range("J34").Value = 1788.6
$ws.Range("L34").Value = 5365.799999999999
$ws.Range("N34").Value = -5533.799999999999
$ws.Range("H39").Value = 117883.555
$ws.Range("J39").Value = 7707.5713
$ws.Range("L39").Value = 23122.7139
$ws.Range("N39").Value = -23710.7139
$ws.Range("H55").Value = 8403425
$ws.Range("I55").Value = 144685.58
$ws.Range("J55").Value = 15629823
$ws.Range("K55").Value = 434056.74
$ws.Range("L55").Value = 46889469
$ws.Range("M55").Value = -433879.74
$ws.Range("N55").Value = -46889823
$ws.Range("H68").Value = 2012.375
$ws.Range("I68").Value = 1666.6666
$ws.Range("J68").Value = 2219.8
$ws.Range("K68").Value = 4999.9998
$ws.Range("L68").Value = 6659.400000000001
$ws.Range("M68").Value = -4188.9998
$ws.Range("N68").Value = -8281.400000000001
$ws.Range("H71").Value = 2012.375
$ws.Range("I71").Value = 1666.6666
$ws.Range("J71").Value = 2219.8
$ws.Range("K71").Value = 14999.9994
$ws.Range("L71").Value = 19978.2
$ws.Range("M71").Value = -10943.9994
$ws.Range("N71").Value = -28090.2
$ws.Range("H120").Value = 10000
$ws.Range("I120").Value = 10000
$ws.Range("K120").Value = 30000
$ws.Range("M120").Value = -25162

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 15833.333
$ws.Range("I57").Value = 15833.333
$ws.Range("K57").Value = 15833.333
$ws.Range("M57").Value = -15013.333
$ws.Range("H58").Value = 23427
$ws.Range("J58").Value = 23998.166
$ws.Range("L58").Value = 23998.166
$ws.Range("N58").Value = -24552.166
$ws.Range("H122").Value = 2588.8696
$ws.Range("I122").Value = 2199.5806
$ws.Range("J122").Value = 3393.4
$ws.Range("K122").Value = 6598.7418
$ws.Range("L122").Value = 10180.2
$ws.Range("M122").Value = -4148.7418
$ws.Range("N122").Value = -15080.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17270.285
$ws.Range("I7").Value = 20261
$ws.Range("K7").Value = 20261
$ws.Range("M7").Value = -20149
$ws.Range("H22").Value = 1477.1666
$ws.Range("J22").Value = 1359
$ws.Range("L22").Value = 1359
$ws.Range("N22").Value = -1949
$ws.Range("H27").Value = 1477.1666
$ws.Range("J27").Value = 1359
$ws.Range("L27").Value = 1359
$ws.Range("N27").Value = -1573
$ws.Range("H46").Value = 16736.549
$ws.Range("J46").Value = 3595.5833
$ws.Range("L46").Value = 3595.5833
$ws.Range("N46").Value = -3971.5833
$ws.Range("H55").Value = 2584.4
$ws.Range("I55").Value = 1421.4
$ws.Range("K55").Value = 1421.4
$ws.Range("M55").Value = -1248.4
$ws.Range("H126").Value = 17270.285
$ws.Range("I126").Value = 20261
$ws.Range("K126").Value = 60783
$ws.Range("M126").Value = -58313

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 41467.5
$ws.Range("I2").Value = 41677.145
$ws.Range("K2").Value = 41677.145
$ws.Range("M2").Value = -41565.145
$ws.Range("H5").Value = 15000000
$ws.Range("I5").Value = 15000000
$ws.Range("K5").Value = 15000000
$ws.Range("M5").Value = -14999888
$ws.Range("H107").Value = 22732938
$ws.Range("I107").Value = 8395.214
$ws.Range("K107").Value = 25185.642
$ws.Range("M107").Value = -23265.642

Write-Output "Applied Leviathan_Profits value updates across all sheets."
